# Scheduled market-data refresh: update Leve profit-calculator sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with latest Universalis
# market-board price snapshots (currentAveragePrice*, Leve prices/profits).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 958.3333
$ws.Range("I19").Value = 457.25
$ws.Range("K19").Value = 457.25
$ws.Range("M19").Value = -282.25
$ws.Range("H107").Value = 712.6923
$ws.Range("I107").Value = 366.5
$ws.Range("K107").Value = 366.5
$ws.Range("M107").Value = 1553.5
$ws.Range("H112").Value = 1983.2084
$ws.Range("J112").Value = 1982.4783
$ws.Range("L112").Value = 5947.4349
$ws.Range("N112").Value = -8163.4349
$ws.Range("H116").Value = 17250
$ws.Range("I116").Value = 51500
$ws.Range("J116").Value = 5833.3335
$ws.Range("K116").Value = 51500
$ws.Range("L116").Value = 5833.3335
$ws.Range("M116").Value = -48058
$ws.Range("N116").Value = -12717.3335
$ws.Range("H129").Value = 1125.1578
$ws.Range("J129").Value = 1172.7059
$ws.Range("L129").Value = 3518.1177
$ws.Range("N129").Value = -13518.1177
$ws.Range("H132").Value = 1262.4694
$ws.Range("I132").Value = 1108.2444
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 3324.7332
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -794.7332000000001
$ws.Range("N132").Value = -14052.5
$ws.Range("H138").Value = 3434.3333
$ws.Range("J138").Value = 3451.9473
$ws.Range("L138").Value = 10355.8419
$ws.Range("N138").Value = -20635.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 927580.2
$ws.Range("J2").Value = 1981.5
$ws.Range("L2").Value = 1981.5
$ws.Range("N2").Value = -2207.5
$ws.Range("H32").Value = 4012.8794
$ws.Range("I32").Value = 3195.2307
$ws.Range("K32").Value = 3195.2307
$ws.Range("M32").Value = -2908.2307
$ws.Range("H102").Value = 1499.25
$ws.Range("I102").Value = 1499.25
$ws.Range("K102").Value = 1499.25
$ws.Range("M102").Value = 122.75
$ws.Range("H116").Value = 927580.2
$ws.Range("J116").Value = 1981.5
$ws.Range("L116").Value = 1981.5
$ws.Range("N116").Value = -6569.5
$ws.Range("H132").Value = 1263.5946
$ws.Range("I132").Value = 1018.28125
$ws.Range("K132").Value = 3054.84375
$ws.Range("M132").Value = -524.84375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 927580.2
$ws.Range("J3").Value = 1981.5
$ws.Range("L3").Value = 1981.5
$ws.Range("N3").Value = -2209.5
$ws.Range("H80").Value = 10680.6
$ws.Range("I80").Value = 450
$ws.Range("J80").Value = 13238.25
$ws.Range("K80").Value = 450
$ws.Range("L80").Value = 13238.25
$ws.Range("M80").Value = 548
$ws.Range("N80").Value = -15234.25
$ws.Range("H83").Value = 10680.6
$ws.Range("I83").Value = 450
$ws.Range("J83").Value = 13238.25
$ws.Range("K83").Value = 2250
$ws.Range("L83").Value = 66191.25
$ws.Range("M83").Value = 2742
$ws.Range("N83").Value = -76175.25
$ws.Range("H105").Value = 2088.5
$ws.Range("I105").Value = 2074.4
$ws.Range("K105").Value = 2074.4
$ws.Range("M105").Value = -327.4000000000001
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2623
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2789
$ws.Range("H134").Value = 1796.1936
$ws.Range("I134").Value = 1658.6296
$ws.Range("J134").Value = 2724.75
$ws.Range("K134").Value = 4975.8888
$ws.Range("L134").Value = 8174.25
$ws.Range("M134").Value = -2440.8888
$ws.Range("N134").Value = -13244.25
$ws.Range("H141").Value = 70996.5
$ws.Range("J141").Value = 69993
$ws.Range("L141").Value = 69993
$ws.Range("N141").Value = -80353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9811.194
$ws.Range("J131").Value = 10309.616
$ws.Range("L131").Value = 30928.848
$ws.Range("N131").Value = -41008.848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20970
$ws.Range("H123").Value = 13240.1
$ws.Range("J123").Value = 13240.1
$ws.Range("L123").Value = 13240.1
$ws.Range("N123").Value = -18140.1
$ws.Range("H132").Value = 1284232.1
$ws.Range("I132").Value = 1833122.5
$ws.Range("J132").Value = 3487.7778
$ws.Range("K132").Value = 5499367.5
$ws.Range("L132").Value = 10463.3334
$ws.Range("M132").Value = -5496837.5
$ws.Range("N132").Value = -15523.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 40999
$ws.Range("J6").Value = 40999
$ws.Range("L6").Value = 40999
$ws.Range("N6").Value = -41223
$ws.Range("H55").Value = 686.55554
$ws.Range("I55").Value = 649.8333
$ws.Range("K55").Value = 649.8333
$ws.Range("M55").Value = -476.8333
$ws.Range("H61").Value = 3589.8333
$ws.Range("I61").Value = 2884.75
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2884.75
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2682.75
$ws.Range("N61").Value = -5404
$ws.Range("H100").Value = 1550
$ws.Range("I100").Value = 1437.5
$ws.Range("K100").Value = 1437.5
$ws.Range("M100").Value = -896.5
$ws.Range("H113").Value = 3589.8333
$ws.Range("I113").Value = 2884.75
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2884.75
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -714.75
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 4641.952
$ws.Range("I136").Value = 3700
$ws.Range("K136").Value = 11100
$ws.Range("M136").Value = -8550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 712.44446
$ws.Range("I100").Value = 516
$ws.Range("K100").Value = 1032
$ws.Range("M100").Value = -491
$ws.Range("H132").Value = 1699.9697
$ws.Range("I132").Value = 960.6957
$ws.Range("K132").Value = 2882.0871
$ws.Range("M132").Value = -352.0870999999997
$ws.Range("H136").Value = 11112418
$ws.Range("I136").Value = 18519874
$ws.Range("J136").Value = 1235.15
$ws.Range("K136").Value = 55559622
$ws.Range("L136").Value = 3705.45
$ws.Range("M136").Value = -55557072
$ws.Range("N136").Value = -8805.450000000001
